$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: row 7 (bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0 handback report) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.cb9fba561be009cb8ae29fa4721a900779ee154f.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-24 10:56:54"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a03d7fb90e7bbc30f1f460e568d25b941bb5611/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e47e919eda2a5a3f259b5f0106c7b95bbdb6d112/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md."

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e47e919eda2a5a3f259b5f0106c7b95bbdb6d112/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md",
    "",
    "",
    "bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md"
) | Out-Null

# --- de-de sheet: row 7 (bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0 handback report) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.cb9fba561be009cb8ae29fa4721a900779ee154f.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-24 10:57:03"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a03d7fb90e7bbc30f1f460e568d25b941bb5611/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e47e919eda2a5a3f259b5f0106c7b95bbdb6d112/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md."

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e47e919eda2a5a3f259b5f0106c7b95bbdb6d112/e2e/bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md",
    "",
    "",
    "bdb6791c-aa83-46b8-8093-a1fd5fdb8ce0.md"
) | Out-Null
